$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "20.241.92"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.50%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.441.89"
$ws.Range("D3").Style = "Normal"

$ws.Range("E4").Value = "  +0.69%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9194"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -8.09%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "274.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.56%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3639"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3079"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.15%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "38.77"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.47%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.017"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.68%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06480"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.03%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9989"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.17%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.315"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.44%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "17.39"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.73%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.017"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.38%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001005"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.57%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.441.75"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.42%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9364"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -6.41%  "

$ws.Range("E19").Value = "  -0.93%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "67.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.21%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.324"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.61%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.15"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.73%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.70"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.26%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.242"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.67%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "20.271.46"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.57%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "139.43"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.25%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.034"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -9.57%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.87"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.19%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.594.46"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.94%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "110.10"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.86%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.017"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.45%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.805"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -9.77%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7793"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.26%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.07658"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.27%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.454"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.56%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05740"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.51%  "

$ws.Range("E37").Value = "  +4.41%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.628"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.28%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01980"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.19%  "

$ws.Range("B40").Value = "Frax"
$ws.Range("C40").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9323"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.85%  "

$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "10.11"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.94%  "

$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1834"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.75%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.928"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -17.52%  "

$ws.Range("B44").Value = "PancakeSwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.479"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.97%  "

$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5174"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "11.80"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.69%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "115.11"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.62%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5080"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.03%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.725"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.48%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06360"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.87%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9880"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.26%  "
